$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A33").Value = "Thomas Cavagna"
$ws.Range("B33").Value = "Thomas Debiasi | Mai una gioia"
$ws.Range("C33").Value = "Thomas Cavagna | Mai una gioia"
$ws.Range("D33").Value = "Luca Frasca | Clitoriders"
$ws.Range("E33").Value = "Federico Nicolodi | U.SGUARNA"
$ws.Range("F33").Value = "Davide  Bazzano | IMONTAGNA"
